# Commit: "PID to PD on presenation"
# Change every visible "PID" (controller) reference to "PD" on the two
# slides that mention the PID/PD controller.

$p = $ppt.ActivePresentation

# --- Slide "Project overview" (2nd slide) ---------------------------------
# Paragraph: "To implement a PID controller to manipulate movement of
# Turtlebot to goal." -> replace PID with PD. The sentence lives in a
# single run, so retargeting the whole run's Text keeps it as one run.
$slideOverview = $p.Slides.Item(2)
$overviewBody = $slideOverview.Shapes.Item(2).TextFrame.TextRange
for ($i = 1; $i -le $overviewBody.Paragraphs().Count; $i++) {
    $para = $overviewBody.Paragraphs($i, 1)
    if ($para.Text -like "*PID controller*") {
        $run = $para.Runs(1, 1)
        $run.Text = $run.Text -replace "PID controller", "PD controller"
    }
}

# --- Slide "Methodology" (4th slide) ---------------------------------------
# Paragraph: "PID Controller " -> "PD Controller ". Only the "PID " prefix
# is retyped, leaving "Controller " as a separate, untouched run.
$slideMethod = $p.Slides.Item(4)
$methodBody = $slideMethod.Shapes.Item(2).TextFrame.TextRange
for ($i = 1; $i -le $methodBody.Paragraphs().Count; $i++) {
    $para = $methodBody.Paragraphs($i, 1)
    if ($para.Text.TrimEnd() -eq "PID Controller") {
        $prefix = $para.Characters(1, 4)
        $prefix.Text = "PD "
    }
}
